$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.058.30"
$ws.Range("E2").Value = "  -4.01%  "
$ws.Range("D3").Value = "1.966.24"
$ws.Range("E3").Value = "  -6.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.63"
$ws.Range("E5").Value = "  -4.08%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4994"
$ws.Range("E7").Value = "  -5.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4218"
$ws.Range("E8").Value = "  -3.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.17"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09212"
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.101"
$ws.Range("E11").Value = "  -6.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.92"
$ws.Range("E12").Value = "  -7.27%  "
$ws.Range("D13").Value = "1.976.78"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.881"
$ws.Range("E14").Value = "  -8.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.452"
$ws.Range("E15").Value = "  -6.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  -4.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.50"
$ws.Range("E18").Value = "  -9.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06700"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.25"
$ws.Range("E20").Value = "  -8.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.961"
$ws.Range("E22").Value = "  -6.17%  "
$ws.Range("D23").Value = "29.086.14"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.08"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.281"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Value = "2.203.01"
$ws.Range("E26").Value = "  -3.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.64"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.63"
$ws.Range("E28").Value = "  -5.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.222"
$ws.Range("E29").Value = "  -9.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.267"
$ws.Range("E30").Value = "  -9.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.97"
$ws.Range("E31").Value = "  -5.10%  "
$ws.Range("E32").Value = "  -7.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09861"
$ws.Range("E33").Value = "  -6.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.533"
$ws.Range("E34").Value = "  -8.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.785"
$ws.Range("E35").Value = "  -7.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.672"
$ws.Range("E36").Value = "  -6.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02426"
$ws.Range("E37").Value = "  -7.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.302"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.952"
$ws.Range("E39").Value = "  -11.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06322"
$ws.Range("E40").Value = "  -6.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6465"
$ws.Range("E41").Value = "  -6.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.46"
$ws.Range("E42").Value = "  -8.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1991"
$ws.Range("E43").Value = "  -10.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6243"
$ws.Range("E45").Value = "  -7.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.41"
$ws.Range("E46").Value = "  -5.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.185"
$ws.Range("E47").Value = "  -8.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.285"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.466"
$ws.Range("E49").Value = "  -4.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000332"
$ws.Range("E50").Value = "  -4.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06918"
$ws.Range("E51").Value = "  -4.75%  "
